$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.711.25"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.022.47"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.13"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("E10").Value = "  +2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.369"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.545.00"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.94%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.30"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.707.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.024.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.81%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +4.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.85"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.152.69"
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0920"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.37"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.66"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.23"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +15.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.83"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.22"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0663"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.061.52"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  +4.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.657"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.203.56"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.53%  "
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("E47").Value = "  +7.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.931"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.79"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.84"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  +1.74%  "
